# Updating STN and SLS sampling frequency and fixing map legend
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 = Summer Townet Survey (STN); Row 15 = Smelt Larva Survey (SLS).
# Both had Frequency (column G) of "2X monthly" -> now "Every 2 weeks".
$ws.Range("G13").Value = "Every 2 weeks"
$ws.Range("G15").Value = "Every 2 weeks"

# Fix map legend: move the active selection from G2 to G15.
$ws.Range("G15").Select()
